# Auto-generated script to apply cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextCell $ws.Range("D2") "66.077.34"
Set-TextCell $ws.Range("E2") "  -1.73%  "
Set-TextCell $ws.Range("D3") "3.416.62"
Set-TextCell $ws.Range("E3") "  -1.26%  "
Set-TextCell $ws.Range("E4") "  +0.00%  "
Set-TextCell $ws.Range("D5") "582.94"
Set-TextCell $ws.Range("E5") "  -1.73%  "
Set-TextCell $ws.Range("D6") "172.44"
Set-TextCell $ws.Range("E6") "  -4.31%  "
Set-TextCell $ws.Range("E7") "  +0.04%  "
Set-TextCell $ws.Range("D8") "0.589"
Set-TextCell $ws.Range("E8") "  -3.27%  "
Set-TextCell $ws.Range("D9") "3.411.47"
Set-TextCell $ws.Range("E9") "  -1.36%  "
Set-TextCell $ws.Range("D10") "0.130"
Set-TextCell $ws.Range("E10") "  -6.87%  "
Set-TextCell $ws.Range("D11") "6.83"
Set-TextCell $ws.Range("E11") "  -1.66%  "
Set-TextCell $ws.Range("D12") "0.408"
Set-TextCell $ws.Range("E12") "  -5.00%  "
Set-TextCell $ws.Range("D13") "4.006.03"
Set-TextCell $ws.Range("E13") "  -1.17%  "
Set-TextCell $ws.Range("E14") "  -0.66%  "
Set-TextCell $ws.Range("D15") "29.79"
Set-TextCell $ws.Range("E15") "  -6.92%  "
Set-TextCell $ws.Range("D16") "66.093.68"
Set-TextCell $ws.Range("E16") "  -1.65%  "
Set-TextCell $ws.Range("D17") "0.0000170"
Set-TextCell $ws.Range("E17") "  -4.15%  "
Set-TextCell $ws.Range("D18") "3.411.90"
Set-TextCell $ws.Range("E18") "  -1.45%  "
Set-TextCell $ws.Range("D19") "5.87"
Set-TextCell $ws.Range("E19") "  -5.30%  "
Set-TextCell $ws.Range("D20") "13.64"
Set-TextCell $ws.Range("E20") "  -3.33%  "
Set-TextCell $ws.Range("D21") "365.50"
Set-TextCell $ws.Range("E21") "  -6.97%  "
Set-TextCell $ws.Range("D22") "7.65"
Set-TextCell $ws.Range("E22") "  -3.30%  "
Set-TextCell $ws.Range("E23") "  +0.18%  "
Set-TextCell $ws.Range("D24") "5.70"
Set-TextCell $ws.Range("E24") "  -1.25%  "
Set-TextCell $ws.Range("D25") "71.05"
Set-TextCell $ws.Range("E25") "  -0.78%  "
Set-TextCell $ws.Range("D26") "0.524"
Set-TextCell $ws.Range("E26") "  -2.58%  "
Set-TextCell $ws.Range("D27") "0.0000118"
Set-TextCell $ws.Range("E27") "  -2.44%  "
Set-TextCell $ws.Range("D28") "9.58"
Set-TextCell $ws.Range("E28") "  -7.78%  "
Set-TextCell $ws.Range("E29") "  +0.84%  "
Set-TextCell $ws.Range("D30") "0.998"
Set-TextCell $ws.Range("E30") "  -0.10%  "
Set-TextCell $ws.Range("B31") "NEARProtocol"
Set-TextCell $ws.Range("C31") "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell $ws.Range("D31") "5.76"
Set-TextCell $ws.Range("E31") "  -5.68%  "
Set-TextCell $ws.Range("B32") "EthereumClassic"
Set-TextCell $ws.Range("C32") "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell $ws.Range("D32") "23.64"
Set-TextCell $ws.Range("E32") "  +0.68%  "
Set-TextCell $ws.Range("D33") "1.98"
Set-TextCell $ws.Range("E33") "  -3.66%  "
Set-TextCell $ws.Range("D34") "0.999"
Set-TextCell $ws.Range("E34") "  -0.10%  "
Set-TextCell $ws.Range("D35") "1.29"
Set-TextCell $ws.Range("E35") "  -8.10%  "
Set-TextCell $ws.Range("D36") "7.01"
Set-TextCell $ws.Range("E36") "  -4.19%  "
Set-TextCell $ws.Range("D37") "1.53"
Set-TextCell $ws.Range("E37") "  -3.18%  "
Set-TextCell $ws.Range("D38") "160.42"
Set-TextCell $ws.Range("E38") "  -0.32%  "
Set-TextCell $ws.Range("D39") "29.03"
Set-TextCell $ws.Range("E39") "  +11.21%  "
Set-TextCell $ws.Range("E40") "  -0.27%  "
Set-TextCell $ws.Range("D41") "2.62"
Set-TextCell $ws.Range("E41") "  -6.98%  "
Set-TextCell $ws.Range("D42") "1.75"
Set-TextCell $ws.Range("E42") "  -6.42%  "
Set-TextCell $ws.Range("D43") "2.710.84"
Set-TextCell $ws.Range("E43") "  -1.55%  "
Set-TextCell $ws.Range("D44") "4.38"
Set-TextCell $ws.Range("E44") "  -5.85%  "
Set-TextCell $ws.Range("D45") "6.30"
Set-TextCell $ws.Range("E45") "  -5.95%  "
Set-TextCell $ws.Range("D46") "0.0677"
Set-TextCell $ws.Range("E46") "  -5.83%  "
Set-TextCell $ws.Range("D47") "39.87"
Set-TextCell $ws.Range("E47") "  -3.63%  "
Set-TextCell $ws.Range("D48") "0.0288"
Set-TextCell $ws.Range("E48") "  -3.26%  "
Set-TextCell $ws.Range("D49") "23.94"
Set-TextCell $ws.Range("E49") "  -8.84%  "
Set-TextCell $ws.Range("D50") "305.87"
Set-TextCell $ws.Range("E50") "  -6.12%  "
Set-TextCell $ws.Range("D51") "0.813"
Set-TextCell $ws.Range("E51") "  -3.47%  "
